$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert D68:D73 into a shared formula block (matches si="2") ---
$ws.Range("D68:D73").Formula = "=C68/(24*60)"

# --- New row 74 (四方坪站 / si 2 shared string) ---
$ws.Range("A74").Value = 45937
$ws.Range("B74").Value = "四方坪站"
$ws.Range("C74").Formula = "=16951/126"
$ws.Range("D74").Formula = "=C74/(24*60)"
$ws.Range("E74").Formula = "=10341.38/126"
$ws.Range("F74").Formula = "=3562.56/126"
$ws.Range("G74").Formula = "=10341.38/(16951/60)"
$ws.Range("H74").Value = 3.1428571428571428

# --- New row 75 (高岭站 / si 3 shared string) ---
$ws.Range("A75").Value = 45937
$ws.Range("B75").Value = "高岭站"
$ws.Range("C75").Formula = "=5616/36"
$ws.Range("D75").Formula = "=C75/(24*60)"
$ws.Range("E75").Formula = "=4056.06/36"
$ws.Range("F75").Formula = "=1094.62/36"
$ws.Range("G75").Formula = "=4056.06/(5616/60)"
$ws.Range("H75").Formula = "=143/36"

# --- Update the selection shown when the sheet is reopened ---
$ws.Range("E77").Select()

Write-Host "Edit applied"
